$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/value updates (safe from Excel numeric auto-detection)
$ws.Range('D2').Value = '28.398.14'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '1.571.88'
$ws.Range('E3').Value = '  +0.02%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  -3.58%  '
$ws.Range('E9').Value = '  -1.89%  '
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('E12').Value = '  +1.43%  '
$ws.Range('D13').Value = '1.797.52'
$ws.Range('E13').Value = '  +0.07%  '
$ws.Range('D14').Value = '1.587.13'
$ws.Range('E14').Value = '  +0.89%  '
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '28.399.96'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('E17').Value = '  -1.10%  '
$ws.Range('E18').Value = '  -0.81%  '
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('D21').Value = '0.0₃0684'
$ws.Range('E21').Value = '  -1.15%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('E23').Value = '  +1.69%  '
$ws.Range('E24').Value = '  -1.02%  '
$ws.Range('E25').Value = '  +1.28%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('E28').Value = '  -0.27%  '
$ws.Range('E29').Value = '  -1.24%  '
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('E31').Value = '  +3.92%  '
$ws.Range('E32').Value = '  -2.36%  '
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('D35').Value = '1.380.62'
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('E36').Value = '  +4.62%  '
$ws.Range('E37').Value = '  -1.85%  '
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('E46').Value = '  -4.31%  '
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('E48').Value = '  -6.17%  '
$ws.Range('D49').Value = '1.709.86'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('E50').Value = '  -0.23%  '

# Values that look like plain numbers ("123.45") must be forced to text so Excel
# does not reinterpret them as numeric values; NumberFormat is restored to Normal
# afterward so no stray cell style is introduced.
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '211.92'
$cell.Style = 'Normal'
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '44.62'
$cell.Style = 'Normal'
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '23.67'
$cell.Style = 'Normal'
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '61.62'
$cell.Style = 'Normal'
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '229.76'
$cell.Style = 'Normal'
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '7.39'
$cell.Style = 'Normal'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '3.95'
$cell.Style = 'Normal'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '151.43'
$cell.Style = 'Normal'
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '14.91'
$cell.Style = 'Normal'
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.0482'
$cell.Style = 'Normal'
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '1.06'
$cell.Style = 'Normal'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.521'
$cell.Style = 'Normal'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.786'
$cell.Style = 'Normal'
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.0473'
$cell.Style = 'Normal'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '62.34'
$cell.Style = 'Normal'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '85.30'
$cell.Style = 'Normal'
